$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, pushing existing rows 99:130 down to 100:131.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly price record
# (same market/product context as the surrounding rows).
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44559
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112031
$ws.Range("G99").Value = "Poroto verde"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 300
$ws.Range("K99").Value = 35000
$ws.Range("L99").Value = 35000
$ws.Range("M99").Value = 35000
$ws.Range("N99").Value = "$/saco 25 kilos"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 1400
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"
